$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 88; this shifts the existing rows 88:188 down to 89:189
# (matching the target dimension growing from A1:R188 to A1:R189).
$ws.Rows.Item(88).Insert()

# Populate the newly inserted row 88 with the new record.
$ws.Cells.Item(88, 1).Value = 10
$ws.Cells.Item(88, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(88, 3).Value = "La Araucanía"
$ws.Cells.Item(88, 4).Value = 44546
$ws.Cells.Item(88, 5).Value = 9
$ws.Cells.Item(88, 6).Value = 100112039
$ws.Cells.Item(88, 7).Value = "Ciboulette"
$ws.Cells.Item(88, 8).Value = "Sin especificar"
$ws.Cells.Item(88, 9).Value = "Primera"
$ws.Cells.Item(88, 10).Value = 65
$ws.Cells.Item(88, 11).Value = 7000
$ws.Cells.Item(88, 12).Value = 7000
$ws.Cells.Item(88, 13).Value = 7000
$ws.Cells.Item(88, 14).Value = "$/docena de atados"
$ws.Cells.Item(88, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(88, 16).Value = 2333
$ws.Cells.Item(88, 17).Value = 3
$ws.Cells.Item(88, 18).Value = "Hortaliza"
